$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 61676.668
$ws.Range("J26").Value = 61676.668
$ws.Range("L26").Value = 61676.668
$ws.Range("N26").Value = -62364.668

$ws.Range("H138").Value = 1722.6222
$ws.Range("I138").Value = 1097.6875
$ws.Range("J138").Value = 1857.7433
$ws.Range("K138").Value = 3293.0625
$ws.Range("L138").Value = 5573.2299
$ws.Range("M138").Value = 1846.9375
$ws.Range("N138").Value = -15853.2299

$ws.Range("H141").Value = 3209.6843
$ws.Range("I141").Value = 2919.4443
$ws.Range("J141").Value = 3470.9
$ws.Range("K141").Value = 8758.332900000001
$ws.Range("L141").Value = 10412.7
$ws.Range("M141").Value = -3578.332900000001
$ws.Range("N141").Value = -20772.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 10000
$ws.Range("I8").Value = 10000
$ws.Range("K8").Value = 10000
$ws.Range("M8").Value = -9856

$ws.Range("H9").Value = 12750
$ws.Range("J9").Value = 12750
$ws.Range("L9").Value = 12750
$ws.Range("N9").Value = -13090

$ws.Range("H20").Value = 12750
$ws.Range("J20").Value = 12750
$ws.Range("L20").Value = 12750
$ws.Range("N20").Value = -13290

$ws.Range("H21").Value = 40224.8
$ws.Range("I21").Value = 1100
$ws.Range("J21").Value = 50006
$ws.Range("K21").Value = 1100
$ws.Range("L21").Value = 50006
$ws.Range("M21").Value = -726
$ws.Range("N21").Value = -50754

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H26").Value = 16121.1
$ws.Range("I26").Value = 1600.5714
$ws.Range("J26").Value = 50002.332
$ws.Range("K26").Value = 1600.5714
$ws.Range("L26").Value = 50002.332
$ws.Range("M26").Value = -1270.5714
$ws.Range("N26").Value = -50662.332

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()

$ws.Range("H39").Value = 23578.857
$ws.Range("I39").Value = 3754
$ws.Range("J39").Value = 50012
$ws.Range("K39").Value = 3754
$ws.Range("L39").Value = 50012
$ws.Range("M39").Value = -3234
$ws.Range("N39").Value = -51052

$ws.Range("H45").Value = 2129.2856
$ws.Range("I45").Value = 1676.5
$ws.Range("J45").Value = 2733
$ws.Range("K45").Value = 1676.5
$ws.Range("L45").Value = 2733
$ws.Range("M45").Value = -1299.5
$ws.Range("N45").Value = -3487

$ws.Range("H61").Value = 3735
$ws.Range("I61").Value = 2477.5
$ws.Range("J61").Value = 6250
$ws.Range("K61").Value = 2477.5
$ws.Range("L61").Value = 6250
$ws.Range("M61").Value = -2265.5
$ws.Range("N61").Value = -6674

$ws.Range("H74").Value = 70115.27
$ws.Range("I74").Value = 75234.96000000001
$ws.Range("J74").Value = 999.5
$ws.Range("K74").Value = 75234.96000000001
$ws.Range("L74").Value = 999.5
$ws.Range("M74").Value = -74360.96000000001
$ws.Range("N74").Value = -2747.5

$ws.Range("H77").Value = 70115.27
$ws.Range("I77").Value = 75234.96000000001
$ws.Range("J77").Value = 999.5
$ws.Range("K77").Value = 376174.8
$ws.Range("L77").Value = 4997.5
$ws.Range("M77").Value = -371806.8
$ws.Range("N77").Value = -13733.5

$ws.Range("H136").Value = 3735
$ws.Range("I136").Value = 2477.5
$ws.Range("J136").Value = 6250
$ws.Range("K136").Value = 7432.5
$ws.Range("L136").Value = 18750
$ws.Range("M136").Value = -4882.5
$ws.Range("N136").Value = -23850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 50000
$ws.Range("J38").Value = 50000
$ws.Range("L38").Value = 50000
$ws.Range("N38").Value = -50832

$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 330.125
$ws.Range("I22").Value = 298.7143
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 298.7143
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = 51.28570000000002
$ws.Range("N22").Value = -1250

$ws.Range("H58").Value = 3018.5
$ws.Range("I58").Value = 953.6818
$ws.Range("J58").Value = 14375
$ws.Range("K58").Value = 953.6818
$ws.Range("L58").Value = 14375
$ws.Range("M58").Value = -750.6818
$ws.Range("N58").Value = -14781

$ws.Range("H107").Value = 248.65218
$ws.Range("I107").Value = 232.6875
$ws.Range("J107").Value = 285.14285
$ws.Range("K107").Value = 232.6875
$ws.Range("L107").Value = 285.14285
$ws.Range("M107").Value = 1687.3125
$ws.Range("N107").Value = -4125.14285

$ws.Range("H136").Value = 3018.5
$ws.Range("I136").Value = 953.6818
$ws.Range("J136").Value = 14375
$ws.Range("K136").Value = 2861.0454
$ws.Range("L136").Value = 43125
$ws.Range("M136").Value = -311.0454
$ws.Range("N136").Value = -48225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 748.6
$ws.Range("J122").Value = 1399.5
$ws.Range("L122").Value = 12595.5
$ws.Range("N122").Value = -17495.5

$ws.Range("H131").Value = 60898216
$ws.Range("J131").Value = 79167544
$ws.Range("L131").Value = 237502632
$ws.Range("N131").Value = -237512712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 56675.168
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 56675.168
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 56675.168
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -57271.168

$ws.Range("H70").Value = 4119.0835
$ws.Range("I70").Value = 4052.8572
$ws.Range("K70").Value = 4052.8572
$ws.Range("M70").Value = -3782.8572

$ws.Range("H73").Value = 4119.0835
$ws.Range("I73").Value = 4052.8572
$ws.Range("K73").Value = 4052.8572
$ws.Range("M73").Value = -3116.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 24850
$ws.Range("J34").Value = 22685.715
$ws.Range("L34").Value = 22685.715
$ws.Range("N34").Value = -23029.715

$ws.Range("H82").Value = 2458.111
$ws.Range("I82").Value = 2900
$ws.Range("J82").Value = 2402.875
$ws.Range("K82").Value = 2900
$ws.Range("L82").Value = 2402.875
$ws.Range("M82").Value = -2539
$ws.Range("N82").Value = -3124.875

$ws.Range("H85").Value = 2458.111
$ws.Range("I85").Value = 2900
$ws.Range("J85").Value = 2402.875
$ws.Range("K85").Value = 2900
$ws.Range("L85").Value = 2402.875
$ws.Range("M85").Value = -1652
$ws.Range("N85").Value = -4898.875

$ws.Range("H122").Value = 2522.2273
$ws.Range("I122").Value = 2242.6155
$ws.Range("J122").Value = 2926.111
$ws.Range("K122").Value = 6727.8465
$ws.Range("L122").Value = 8778.332999999999
$ws.Range("M122").Value = -4277.8465
$ws.Range("N122").Value = -13678.333

$ws.Range("H132").Value = 752981.75
$ws.Range("I132").Value = 202376.5
$ws.Range("J132").Value = 1670657.1
$ws.Range("K132").Value = 607129.5
$ws.Range("L132").Value = 5011971.300000001
$ws.Range("M132").Value = -604599.5
$ws.Range("N132").Value = -5017031.300000001
